$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3356204611722832
$ws.Range("C2").Value = 0.592858919452767
$ws.Range("D2").Value = 0.4584924415561159
$ws.Range("E2").Value = 0.6771206994001261
$ws.Range("F2").Value = 0.6051409294670009
$ws.Range("G2").Value = 18

$ws.Range("B3").Value = 0.1712146867552044
$ws.Range("C3").Value = 0.4703347406954927
$ws.Range("D3").Value = 0.3274315915576724
$ws.Range("E3").Value = 0.572216385257948
$ws.Range("F3").Value = 0.5628049775537718
$ws.Range("G3").Value = 17

$ws.Range("B4").Value = 0.2052142876843926
$ws.Range("C4").Value = 0.4129592675069453
$ws.Range("D4").Value = 0.2759880906773418
$ws.Range("E4").Value = 0.5253456868361458
$ws.Range("F4").Value = 0.4994665813926205
$ws.Range("G4").Value = 16

$ws.Range("B5").Value = 0.3686756894723587
$ws.Range("C5").Value = 0.3990438243463547
$ws.Range("D5").Value = 0.2123069440861257
$ws.Range("E5").Value = 0.4607677767445611
$ws.Range("F5").Value = 0.2860791225684026
$ws.Range("G5").Value = 15

$ws.Range("B6").Value = 0.341617976538504
$ws.Range("C6").Value = 0.3737653469981676
$ws.Range("D6").Value = 0.2000610633735307
$ws.Range("E6").Value = 0.4472818612167619
$ws.Range("F6").Value = 0.2996170762194463
$ws.Range("G6").Value = 14

$ws.Range("B7").Value = 0.3537113215538176
$ws.Range("C7").Value = 0.3951185885610565
$ws.Range("D7").Value = 0.2147690549179593
$ws.Range("E7").Value = 0.4634318233763833
$ws.Range("F7").Value = 0.3116549409996502
$ws.Range("G7").Value = 13

$ws.Range("B8").Value = 0.3810204497723156
$ws.Range("C8").Value = 0.4291397075636965
$ws.Range("D8").Value = 0.2292934380982358
$ws.Range("E8").Value = 0.4788459440135583
$ws.Range("F8").Value = 0.3029254722988744
$ws.Range("G8").Value = 12

$ws.Range("B9").Value = 0.4264878578050817
$ws.Range("C9").Value = 0.4469030973498441
$ws.Range("D9").Value = 0.2370354610461564
$ws.Range("E9").Value = 0.4868628770466653
$ws.Range("F9").Value = 0.2462882965349502
$ws.Range("G9").Value = 11

$ws.Range("B10").Value = 0.3953124032861549
$ws.Range("C10").Value = 0.4125113567458484
$ws.Range("D10").Value = 0.2047118874619409
$ws.Range("E10").Value = 0.4524509779655039
$ws.Range("F10").Value = 0.2319961476453755
$ws.Range("G10").Value = 10

$ws.Range("B11").Value = 0.3698904567253822
$ws.Range("C11").Value = 0.3826200248740704
$ws.Range("D11").Value = 0.1850206325475525
$ws.Range("E11").Value = 0.4301402475327699
$ws.Range("F11").Value = 0.2328666848057507
$ws.Range("G11").Value = 9
